# Sync Ekonomi sheet values/labels with the current "source of truth" and
# make the Ekonomi tab the active/selected sheet (createExcel.mjs sync).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Material")
$ws2 = $wb.Worksheets.Item("Ekonomi")

# --- Update the Ekonomi parameter table -----------------------------------
# Row 2: PrisTimmerIn - inköpspris timmer, now per m2 (area) instead of per m (length)
$ws2.Range("A2").Value = "PrisTimmerIn"
$ws2.Range("B2").Value = 458
$ws2.Range("C2").Value = "Inköpspris timmer (kr/m2)"

# Row 3: PrisTimmerUt - försäljningspris timmer, now per m2 (area) instead of per m (length)
$ws2.Range("A3").Value = "PrisTimmerUt"
$ws2.Range("B3").Value = 850
$ws2.Range("C3").Value = "Försäljningspris timmer (kr/m2)"

# Row 4: Timkostnad
$ws2.Range("A4").Value = "Timkostnad"
$ws2.Range("B4").Value = 550
$ws2.Range("C4").Value = "Timkostnad (kr/h)"

# Row 5: MomsPct
$ws2.Range("A5").Value = "MomsPct"
$ws2.Range("B5").Value = 25
$ws2.Range("C5").Value = "Moms (%)"

# Column A best-fits the longer parameter names now ("PrisTimmerIn"/"PrisTimmerUt")
$ws2.Columns("A").ColumnWidth = 11.75

# --- Make "Ekonomi" the active/selected sheet ------------------------------
$ws2.Activate()
$ws2.Range("B2").Select()
$excel.ActiveWindow.Zoom = 150
